$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2a = New-Object 'object[,]' 1,5
$row2a[0,0] = 1.02
$row2a[0,1] = 1.074712426852077
$row2a[0,2] = 1.077364095165504
$row2a[0,3] = 1.078742575687636
$row2a[0,4] = 1.089610755915166
$ws.Range("B2:F2").Value = $row2a

$row2b = New-Object 'object[,]' 1,6
$row2b[0,0] = 1.056799615230411
$row2b[0,1] = 1.079620530089919
$row2b[0,2] = 1.080045589661624
$row2b[0,3] = 1.081420449699227
$row2b[0,4] = 1.092260431566832
$row2b[0,5] = 1.081153715076852
$ws.Range("I2:N2").Value = $row2b

$row3a = New-Object 'object[,]' 1,5
$row3a[0,0] = 1.02
$row3a[0,1] = 1.075967026908906
$row3a[0,2] = 1.078388025474395
$row3a[0,3] = 1.079858876371388
$row3a[0,4] = 1.090800604478367
$ws.Range("B3:F3").Value = $row3a

$row3b = New-Object 'object[,]' 1,6
$row3b[0,0] = 1.057194358293545
$row3b[0,1] = 1.080532872684834
$row3b[0,2] = 1.080886695480708
$row3b[0,3] = 1.082353962973314
$row3b[0,4] = 1.093269359788487
$row3b[0,5] = 1.082067353302899
$ws.Range("I3:N3").Value = $row3b

$row4a = New-Object 'object[,]' 1,5
$row4a[0,0] = 1.02
$row4a[0,1] = 1.07677866800374
$row4a[0,2] = 1.07905038443934
$row4a[0,3] = 1.080581290846343
$row4a[0,4] = 1.09157074128049
$ws.Range("B4:F4").Value = $row4a

$row4b = New-Object 'object[,]' 1,6
$row4b[0,0] = 1.057448470195458
$row4b[0,1] = 1.081122492562498
$row4b[0,2] = 1.081430142030156
$row4b[0,3] = 1.082957506738167
$row4b[0,4] = 1.093921845742922
$row4b[0,5] = 1.08265781050842
$ws.Range("I4:N4").Value = $row4b

$row5a = New-Object 'object[,]' 1,5
$row5a[0,0] = 1.02
$row5a[0,1] = 1.077119843488723
$row5a[0,2] = 1.079328795084145
$row5a[0,3] = 1.080885017378887
$row5a[0,4] = 1.091894561958632
$ws.Range("B5:F5").Value = $row5a

$row5b = New-Object 'object[,]' 1,6
$row5b[0,0] = 1.057554985003166
$row5b[0,1] = 1.081370195971939
$row5b[0,2] = 1.081658415226881
$row5b[0,3] = 1.083211117522797
$row5b[0,4] = 1.094196066265111
$row5b[0,5] = 1.082905865685114
$ws.Range("I5:N5").Value = $row5b

$row6a = New-Object 'object[,]' 1,5
$row6a[0,0] = 1.02
$row6a[0,1] = 1.077177126159492
$row6a[0,2] = 1.079375538842962
$row6a[0,3] = 1.08093601581102
$row6a[0,4] = 1.091948936146675
$ws.Range("B6:F6").Value = $row6a

$row6b = New-Object 'object[,]' 1,6
$row6b[0,0] = 1.057572850920885
$row6b[0,1] = 1.081411776367119
$row6b[0,2] = 1.081696732087537
$row6b[0,3] = 1.083253692942557
$row6b[0,4] = 1.094242104171439
$row6b[0,5] = 1.082947505129225
$ws.Range("I6:N6").Value = $row6b

$row7a = New-Object 'object[,]' 1,5
$row7a[0,0] = 1.02
$row7a[0,1] = 1.076783226955723
$row7a[0,2] = 1.079054104751822
$row7a[0,3] = 1.080585349161402
$row7a[0,4] = 1.091575067969907
$ws.Range("B7:F7").Value = $row7a

$row7b = New-Object 'object[,]' 1,6
$row7b[0,0] = 1.05744989468399
$row7b[0,1] = 1.081125803065004
$row7b[0,2] = 1.081433192979881
$row7b[0,3] = 1.082960895964146
$row7b[0,4] = 1.093925510221836
$row7b[0,5] = 1.082661125712219
$ws.Range("I7:N7").Value = $row7b

$row8a = New-Object 'object[,]' 1,5
$row8a[0,0] = 1.02
$row8a[0,1] = 1.07513646104021
$row8a[0,2] = 1.077710177650358
$row8a[0,3] = 1.079119816109415
$row8a[0,4] = 1.090012824926839
$ws.Range("B8:F8").Value = $row8a

$row8b = New-Object 'object[,]' 1,6
$row8b[0,0] = 1.05693329275218
$row8b[0,1] = 1.079929011754597
$row8b[0,2] = 1.080330012641268
$row8b[0,3] = 1.081736039435972
$row8b[0,4] = 1.092601478416198
$row8b[0,5] = 1.081462634820877
$ws.Range("I8:N8").Value = $row8b

$row9a = New-Object 'object[,]' 1,5
$row9a[0,0] = 1.02
$row9a[0,1] = 1.072233263187137
$row9a[0,2] = 1.075340487533487
$row9a[0,3] = 1.076538022041547
$row9a[0,4] = 1.087261621701978
$ws.Range("B9:F9").Value = $row9a

$row9b = New-Object 'object[,]' 1,6
$row9b[0,0] = 1.056012890924281
$row9b[0,1] = 1.077814491880857
$row9b[0,2] = 1.078379855217641
$row9b[0,3] = 1.07957379317283
$row9b[0,4] = 1.090265570998644
$row9b[0,5] = 1.079345112086378
$ws.Range("I9:N9").Value = $row9b

$row10a = New-Object 'object[,]' 1,5
$row10a[0,0] = 1.02
$row10a[0,1] = 1.070296736325212
$row10a[0,2] = 1.073759607170955
$row10a[0,3] = 1.074817196021593
$row10a[0,4] = 1.085428526766643
$ws.Range("B10:F10").Value = $row10a

$row10b = New-Object 'object[,]' 1,6
$row10b[0,0] = 1.055392473753653
$row10b[0,1] = 1.076400954802082
$row10b[0,2] = 1.077075506502037
$row10b[0,3] = 1.078129606667551
$row10b[0,4] = 1.08870634588314
$row10b[0,5] = 1.077929567622816
$ws.Range("I10:N10").Value = $row10b

$row11a = New-Object 'object[,]' 1,5
$row11a[0,0] = 1.02
$row11a[0,1] = 1.069457917436599
$row11a[0,2] = 1.073074795271289
$row11a[0,3] = 1.074072129550726
$row11a[0,4] = 1.084635004298445
$ws.Range("B11:F11").Value = $row11a

$row11b = New-Object 'object[,]' 1,6
$row11b[0,0] = 1.055122201066858
$row11b[0,1] = 1.075787944701052
$row11b[0,2] = 1.076509687021795
$row11b[0,3] = 1.07750360347333
$row11b[0,4] = 1.088030703442758
$row11b[0,5] = 1.077315686977129
$ws.Range("I11:N11").Value = $row11b

$row12a = New-Object 'object[,]' 1,5
$row12a[0,0] = 1.02
$row12a[0,1] = 1.069146296544519
$row12a[0,2] = 1.072820382103448
$row12a[0,3] = 1.073795386475441
$row12a[0,4] = 1.084340286004332
$ws.Range("B12:F12").Value = $row12a

$row12b = New-Object 'object[,]' 1,6
$row12b[0,0] = 1.055021564372757
$row12b[0,1] = 1.075560102697668
$row12b[0,2] = 1.076299360616543
$row12b[0,3] = 1.077270977188441
$row12b[0,4] = 1.087779664893006
$row12b[0,5] = 1.077087521411967
$ws.Range("I12:N12").Value = $row12b

$row13a = New-Object 'object[,]' 1,5
$row13a[0,0] = 1.02
$row13a[0,1] = 1.069213142417196
$row13a[0,2] = 1.072874956584723
$row13a[0,3] = 1.073854748478106
$row13a[0,4] = 1.084403502691118
$ws.Range("B13:F13").Value = $row13a

$row13b = New-Object 'object[,]' 1,6
$row13b[0,0] = 1.055043162402306
$row13b[0,1] = 1.075608982067445
$row13b[0,2] = 1.076344483420301
$row13b[0,3] = 1.077320880903432
$row13b[0,4] = 1.087833516927548
$row13b[0,5] = 1.077136470196056
$ws.Range("I13:N13").Value = $row13b

$row14a = New-Object 'object[,]' 1,5
$row14a[0,0] = 1.02
$row14a[0,1] = 1.069432159714384
$row14a[0,2] = 1.073053766295898
$row14a[0,3] = 1.074049253732296
$row14a[0,4] = 1.084610642164195
$ws.Range("B14:F14").Value = $row14a

$row14b = New-Object 'object[,]' 1,6
$row14b[0,0] = 1.055113887416737
$row14b[0,1] = 1.075769114131747
$row14b[0,2] = 1.076492304563247
$row14b[0,3] = 1.07748437657774
$row14b[0,4] = 1.088009954054383
$row14b[0,5] = 1.077296829666256
$ws.Range("I14:N14").Value = $row14b

$row15a = New-Object 'object[,]' 1,5
$row15a[0,0] = 1.02
$row15a[0,1] = 1.069567097297494
$row15a[0,2] = 1.0731639310485
$row15a[0,3] = 1.074169095822339
$row15a[0,4] = 1.084738271723375
$ws.Range("B15:F15").Value = $row15a

$row15b = New-Object 'object[,]' 1,6
$row15b[0,0] = 1.055157430892542
$row15b[0,1] = 1.075867757815393
$row15b[0,2] = 1.076583361352715
$row15b[0,3] = 1.077585098252447
$row15b[0,4] = 1.088118652823304
$row15b[0,5] = 1.077395613435248
$ws.Range("I15:N15").Value = $row15b

$row16a = New-Object 'object[,]' 1,5
$row16a[0,0] = 1.02
$row16a[0,1] = 1.070352399637288
$row16a[0,2] = 1.07380504985446
$row16a[0,3] = 1.074866644809874
$row16a[0,4] = 1.085481194721406
$ws.Range("B16:F16").Value = $row16a

$row16b = New-Object 'object[,]' 1,6
$row16b[0,0] = 1.055410376494104
$row16b[0,1] = 1.076441618334517
$row16b[0,2] = 1.077113036299241
$row16b[0,3] = 1.07817113838262
$row16b[0,4] = 1.088751175635761
$row16b[0,5] = 1.077970288902132
$ws.Range("I16:N16").Value = $row16b

$row17a = New-Object 'object[,]' 1,5
$row17a[0,0] = 1.02
$row17a[0,1] = 1.070844919547148
$row17a[0,2] = 1.074207130698447
$row17a[0,3] = 1.075304214470534
$row17a[0,4] = 1.085947268298621
$ws.Range("B17:F17").Value = $row17a

$row17b = New-Object 'object[,]' 1,6
$row17b[0,0] = 1.055568606170327
$row17b[0,1] = 1.076801333256728
$row17b[0,2] = 1.077445011197876
$row17b[0,3] = 1.078538567848663
$row17b[0,4] = 1.089147808793633
$row17b[0,5] = 1.078330514660792
$ws.Range("I17:N17").Value = $row17b

$row18a = New-Object 'object[,]' 1,5
$row18a[0,0] = 1.02
$row18a[0,1] = 1.071132170278944
$row18a[0,2] = 1.074441630707994
$row18a[0,3] = 1.075559447596099
$row18a[0,4] = 1.086219142450525
$ws.Range("B18:F18").Value = $row18a

$row18b = New-Object 'object[,]' 1,6
$row18b[0,0] = 1.055660741843888
$row18b[0,1] = 1.077011058376727
$row18b[0,2] = 1.077638547461998
$row18b[0,3] = 1.078752819576526
$row18b[0,4] = 1.089379111259604
$row18b[0,5] = 1.078540537614517
$ws.Range("I18:N18").Value = $row18b

$row19a = New-Object 'object[,]' 1,5
$row19a[0,0] = 1.02
$row19a[0,1] = 1.071230110719499
$row19a[0,2] = 1.074521584623763
$row19a[0,3] = 1.075646476586711
$row19a[0,4] = 1.086311848178438
$ws.Range("B19:F19").Value = $row19a

$row19b = New-Object 'object[,]' 1,6
$row19b[0,0] = 1.055692131118361
$row19b[0,1] = 1.077082553920575
$row19b[0,2] = 1.077704521564163
$row19b[0,3] = 1.078825863126444
$row19b[0,4] = 1.089457971508303
$row19b[0,5] = 1.078612134690238
$ws.Range("I19:N19").Value = $row19b

$row20a = New-Object 'object[,]' 1,5
$row20a[0,0] = 1.02
$row20a[0,1] = 1.070792079717696
$row20a[0,2] = 1.074163994027392
$row20a[0,3] = 1.075257266764073
$row20a[0,4] = 1.085897260823066
$ws.Range("B20:F20").Value = $row20a

$row20b = New-Object 'object[,]' 1,6
$row20b[0,0] = 1.055551645871108
$row20b[0,1] = 1.076762748624448
$row20b[0,2] = 1.07740940369231
$row20b[0,3] = 1.078499152741268
$row20b[0,4] = 1.089105258677707
$row20b[0,5] = 1.078291875233908
$ws.Range("I20:N20").Value = $row20b

$row21a = New-Object 'object[,]' 1,5
$row21a[0,0] = 1.02
$row21a[0,1] = 1.069367665876045
$row21a[0,2] = 1.07300111250601
$row21a[0,3] = 1.073991976580883
$row21a[0,4] = 1.084549643902568
$ws.Range("B21:F21").Value = $row21a

$row21b = New-Object 'object[,]' 1,6
$row21b[0,0] = 1.055093067443015
$row21b[0,1] = 1.075721963184078
$row21b[0,2] = 1.076448779234915
$row21b[0,3] = 1.077436233972284
$row21b[0,4] = 1.087957999797575
$row21b[0,5] = 1.077249611758833
$ws.Range("I21:N21").Value = $row21b

$row22a = New-Object 'object[,]' 1,5
$row22a[0,0] = 1.02
$row22a[0,1] = 1.068471811700241
$row22a[0,2] = 1.072269710136006
$row22a[0,3] = 1.073196482670385
$row22a[0,4] = 1.083702523786943
$ws.Range("B22:F22").Value = $row22a

$row22b = New-Object 'object[,]' 1,6
$row22b[0,0] = 1.054803320868781
$row22b[0,1] = 1.075066753298158
$row22b[0,2] = 1.075843894059746
$row22b[0,3] = 1.076767350647457
$row22b[0,4] = 1.087236238669578
$row22b[0,5] = 1.076593471399721
$ws.Range("I22:N22").Value = $row22b

$row23a = New-Object 'object[,]' 1,5
$row23a[0,0] = 1.02
$row23a[0,1] = 1.068946747298317
$row23a[0,2] = 1.072657464793446
$row23a[0,3] = 1.073618185388931
$row23a[0,4] = 1.084151581693326
$ws.Range("B23:F23").Value = $row23a

$row23b = New-Object 'object[,]' 1,6
$row23b[0,0] = 1.054957055888587
$row23b[0,1] = 1.075414171376215
$row23b[0,2] = 1.076164641156371
$row23b[0,3] = 1.077121994232896
$row23b[0,4] = 1.087618899566041
$row23b[0,5] = 1.076941382851299
$ws.Range("I23:N23").Value = $row23b

$row24a = New-Object 'object[,]' 1,5
$row24a[0,0] = 1.02
$row24a[0,1] = 1.070815955848069
$row24a[0,2] = 1.074183485715598
$row24a[0,3] = 1.07527848039289
$row24a[0,4] = 1.085919856982948
$ws.Range("B24:F24").Value = $row24a

$row24b = New-Object 'object[,]' 1,6
$row24b[0,0] = 1.055559309985477
$row24b[0,1] = 1.076780183638629
$row24b[0,2] = 1.077425493497566
$row24b[0,3] = 1.078516962928475
$row24b[0,4] = 1.089124485388675
$row24b[0,5] = 1.078309335007809
$ws.Range("I24:N24").Value = $row24b

$row25a = New-Object 'object[,]' 1,5
$row25a[0,0] = 1.02
$row25a[0,1] = 1.072983987017914
$row25a[0,2] = 1.075953295979195
$row25a[0,3] = 1.077205406938988
$row25a[0,4] = 1.087972684304234
$ws.Range("B25:F25").Value = $row25a

$row25b = New-Object 'object[,]' 1,6
$row25b[0,0] = 1.056252035535476
$row25b[0,1] = 1.078361820068317
$row25b[0,2] = 1.07888476075272
$row25b[0,3] = 1.080133254338202
$row25b[0,4] = 1.090869797880398
$row25b[0,5] = 1.079893217542644
$ws.Range("I25:N25").Value = $row25b
